# Adding the changes we made on may 9th
#
# The data table is a sliding window of accelerometer samples. The window is
# advanced by 13 rows: 13 brand-new rows of (x, y, z) samples are inserted at
# the top of the data (directly under the header row), the previously-existing
# rows shift down by 13, and the 3 oldest rows that fall off the end of the
# window are removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to insert, in order, starting at row 2 (right after the header row).
$newRows = @(
    @(-3.616065740585328, 4.469793319702149, 0.258730050176382),
    @(-3.719920873641968, 4.70874035358429,  0.04408367723226582),
    @(-3.721616864204406, 4.585918724536896, 0.2726323418319225),
    @(-3.805010795593262, 4.544945240020752, 0.311984956264496),
    @(-3.887511849403381, 4.449418604373932, 0.4409204423427582),
    @(-3.778247833251953, 4.410304188728333, 0.5136718302965164),
    @(-3.584390580654144, 4.580866992473602, 0.3528684750199313),
    @(-3.46419882774353,  4.518833875656127, 0.403674334287644),
    @(-3.414171874523162, 4.371547281742096, 0.5634630396962166),
    @(-3.423850417137146, 4.383899688720703, 0.5505108982324599),
    @(-3.575843572616578, 4.333066165447235, 0.460273951292038),
    @(-3.668661117553711, 4.494052410125732, 0.1578152179718018),
    @(-3.737768590450287, 4.342036247253418, 0.1853629685938359)
)

$insertCount = $newRows.Count

# Insert blank rows before row 2, shifting existing data (rows 2-21) down by $insertCount.
$insertRange = $ws.Range("A2:C$($insertCount + 1)")
$insertRange.EntireRow.Insert()

# The inserted rows pick up formatting from the header row above; clear it so the
# new data cells remain unstyled, matching the rest of the plain data rows.
$ws.Range("A2:C$($insertCount + 1)").ClearFormats()

# Fill the newly inserted rows with the new data values.
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowIndex = 2 + $i
    $row = $newRows[$i]
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
}

# The window only keeps 30 rows of data; the 3 oldest rows that were pushed past
# the end (originally rows 19-21, now rows 32-34) are removed entirely.
$ws.Range("A32:C34").EntireRow.Delete()
